# Add the "2022-Q1" holdings sheet (before the "总计" summary sheet) and
# prepend a matching row to the "总计" sheet's history table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet named "2022-Q1" right before the "总计" sheet.
# ---------------------------------------------------------------------
$summarySheetBefore = $wb.Worksheets.Item($wb.Worksheets.Count)   # "总计" (last tab)
$templateSheet = $wb.Worksheets.Item($wb.Worksheets.Count - 1)  # "2021-Q4" (same layout to copy formats from)

$newSheet = $wb.Worksheets.Add($summarySheetBefore)
$newSheet.Name = "2022-Q1"

# NOTE: the COM shim's `Worksheets.Add(Before)` handle aliases whatever
# object reference was passed in as `Before` once the insertion happens, so
# re-resolve the "总计" sheet fresh, by name, rather than reusing the
# pre-insertion reference.
$summarySheet = $wb.Worksheets.Item("总计")

# Copy the header row + row-index column formatting from the "2021-Q4" sheet
# (bold, centered, bordered "s=2" style) so the new sheet matches the others.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$templateSheet.Range("A2:A10").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding data, rows 2-10 (A = row index 0..8, B = fund code, C = fund
# name, D = fund size, E = total stock position, F = position share,
# G = holding value in 100M CNY, H = position rank). B/D/E/F/G are numeric
# looking text in the source data (fund codes keep leading zeros, the rest
# keep fixed decimal formatting) so they are entered with a leading "'" to
# force text the way the original sheets store them; only G10 is a genuine
# 0 number in the source.
$rows = @(
    @{ A=0; B="160314"; C="华夏行业混合(LOF)";                 D="21.89"; E="92.15"; F="7.45"; G="1.6308"; H=1 },
    @{ A=1; B="004693"; C="新疆前海联合泳隽灵活配置混合A";     D="9.08";  E="93.74"; F="4.69"; G="0.4259"; H=7 },
    @{ A=2; B="004640"; C="华夏节能环保股票";                   D="5.45";  E="90.73"; F="4.19"; G="0.2284"; H=8 },
    @{ A=3; B="011160"; C="富国质量成长6个月持有期混合A";       D="3.80";  E="91.55"; F="2.30"; G="0.0874"; H=6 },
    @{ A=4; B="000963"; C="兴业多策略灵活配置混合";             D="2.07";  E="75.34"; F="3.64"; G="0.0753"; H=5 },
    @{ A=5; B="005933"; C="新疆前海联合先进制造灵活配置混合A"; D="1.26";  E="89.79"; F="4.46"; G="0.0562"; H=9 },
    @{ A=6; B="005934"; C="新疆前海联合先进制造灵活配置混合C"; D="0.09";  E="89.79"; F="4.46"; G="0.0040"; H=9 },
    @{ A=7; B="011161"; C="富国质量成长6个月持有期混合C";       D="0.12";  E="91.55"; F="2.30"; G="0.0028"; H=6 },
    @{ A=8; B="007042"; C="新疆前海联合泳隽灵活配置混合C";     D="0.00";  E="93.74"; F="4.69"; G=0;        H=7 }
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row.A
    $newSheet.Range("B$r").Value = "'" + $row.B
    $newSheet.Range("C$r").Value = $row.C
    $newSheet.Range("D$r").Value = "'" + $row.D
    $newSheet.Range("E$r").Value = "'" + $row.E
    $newSheet.Range("F$r").Value = "'" + $row.F
    if ($row.G -is [string]) {
        $newSheet.Range("G$r").Value = "'" + $row.G
    } else {
        $newSheet.Range("G$r").Value = $row.G
    }
    $newSheet.Range("H$r").Value = $row.H
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" (summary) sheet.
# ---------------------------------------------------------------------
$summarySheet.Rows(2).Insert()
$summarySheet.Range("B2:D2").ClearFormats()

# Re-apply the row-index column style ("s=2") to the new A2 cell by
# copying formats from the row directly below (the old A2, now A3).
$summarySheet.Range("A3").Copy()
$summarySheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$summarySheet.Range("A2").Value = 0
$summarySheet.Range("B2").Value = "2022-Q1"
$summarySheet.Range("C2").Value = 9
$summarySheet.Range("D2").Value = 2.51

$excel.CutCopyMode = 0
